$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Table 2 ("Constrained-3" test cases): first column values 1..6
# become 7..12 (index shifted by +6).
# ------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$mapT2 = @{2="7"; 3="8"; 4="9"; 5="10"; 6="11"; 7="12"}
foreach ($row in 2..7) {
    $cell = $t2.Cell($row, 1)
    $rng = $cell.Range
    $target = $d.Range($rng.Start, $rng.Start + 1)
    $target.Text = $mapT2[$row]
}

# ------------------------------------------------------------------
# Table 3 ("Constrained-4" test cases): first column values 1..7
# become 13..19 (index shifted by +12). The "_GoBack" bookmark that
# used to sit in the trailing empty paragraph at the end of the
# document now wraps the last updated cell's text ("19").
# ------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$mapT3 = @{2="13"; 3="14"; 4="15"; 5="16"; 6="17"; 7="18"; 8="19"}
foreach ($row in 2..8) {
    $cell = $t3.Cell($row, 1)
    $rng = $cell.Range
    $target = $d.Range($rng.Start, $rng.Start + 1)

    if ($row -eq 8) {
        # Move the "_GoBack" bookmark onto this run before changing its
        # text so it ends up wrapping the new "19" value.
        $d.Bookmarks.Add("_GoBack", $target)
        $target2 = $d.Range($rng.Start, $rng.Start + 1)
        $target2.Text = $mapT3[$row]
    } else {
        $target.Text = $mapT3[$row]
    }
}
